$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (approximates Excel auto-fit widths after data change)
$ws.Columns.Item(1).ColumnWidth = 25.85546875
$ws.Columns.Item(2).ColumnWidth = 33.28515625
$ws.Columns.Item(3).ColumnWidth = 31.85546875
$ws.Columns.Item(4).ColumnWidth = 38.7109375
$ws.Columns.Item(5).ColumnWidth = 32
$ws.Columns.Item(6).ColumnWidth = 30.5703125
$ws.Columns.Item(7).ColumnWidth = 37.42578125
$ws.Columns.Item(8).ColumnWidth = 32.28515625
$ws.Columns.Item(9).ColumnWidth = 31

# Update data values (hose length sweep results for 0.5s timestep run)
$ws.Range("A2").Value = 5
$ws.Range("B2").Value = 23.719073423116885
$ws.Range("C2").Value = 16.400651092021075
$ws.Range("D2").Value = 0.69145412215119706
$ws.Range("E2").Value = 21.219838883106998
$ws.Range("F2").Value = 15.364052826364389
$ws.Range("G2").Value = 0.72404191714177579
$ws.Range("H2").Value = 221
$ws.Range("I2").Value = 187.5
$ws.Range("A3").Value = 6.666666666666667
$ws.Range("B3").Value = 25.009641843464983
$ws.Range("C3").Value = 16.457429892668891
$ws.Range("D3").Value = 0.65804340564633934
$ws.Range("E3").Value = 22.472736591637485
$ws.Range("F3").Value = 15.440103734246339
$ws.Range("G3").Value = 0.68705934728002294
$ws.Range("H3").Value = 230
$ws.Range("I3").Value = 197
$ws.Range("A4").Value = 8.3333333333333339
$ws.Range("B4").Value = 26.006422696088226
$ws.Range("C4").Value = 16.389952336585843
$ws.Range("D4").Value = 0.6302270992100405
$ws.Range("E4").Value = 23.180308782811931
$ws.Range("F4").Value = 15.188744381416281
$ws.Range("G4").Value = 0.65524340179103435
$ws.Range("H4").Value = 237.5
$ws.Range("I4").Value = 204
$ws.Range("A5").Value = 10
$ws.Range("B5").Value = 27.031550082488675
$ws.Range("C5").Value = 16.34136219020867
$ws.Range("D5").Value = 0.60452923122580304
$ws.Range("E5").Value = 24.269764332926869
$ws.Range("F5").Value = 15.247376752095644
$ws.Range("G5").Value = 0.62824576880664129
$ws.Range("H5").Value = 245
$ws.Range("I5").Value = 211.5
$ws.Range("A6").Value = 11.666666666666668
$ws.Range("B6").Value = 28.088721312609767
$ws.Range("C6").Value = 16.353690660982164
$ws.Range("D6").Value = 0.5822155618611432
$ws.Range("E6").Value = 25.411963687603777
$ws.Range("F6").Value = 15.280473821173853
$ws.Range("G6").Value = 0.60131023359787927
$ws.Range("H6").Value = 252
$ws.Range("I6").Value = 218.5
$ws.Range("A7").Value = 13.333333333333334
$ws.Range("B7").Value = 29.181120840689157
$ws.Range("C7").Value = 16.338879713053931
$ws.Range("D7").Value = 0.55991268472016864
$ws.Range("E7").Value = 26.575267604869079
$ws.Range("F7").Value = 15.399096810506128
$ws.Range("G7").Value = 0.57945218236239804
$ws.Range("H7").Value = 258.5
$ws.Range("I7").Value = 225
$ws.Range("A8").Value = 15
$ws.Range("B8").Value = 30.262837690749702
$ws.Range("C8").Value = 16.376138598143829
$ws.Range("D8").Value = 0.54113030527700468
$ws.Range("E8").Value = 27.457945904237658
$ws.Range("F8").Value = 15.265649653346042
$ws.Range("G8").Value = 0.55596473627657828
$ws.Range("H8").Value = 265
$ws.Range("I8").Value = 230.5
$ws.Range("A9").Value = 16.666666666666664
$ws.Range("B9").Value = 31.023834898998732
$ws.Range("C9").Value = 16.209673096757516
$ws.Range("D9").Value = 0.52249095411736701
$ws.Range("E9").Value = 28.361771424029445
$ws.Range("F9").Value = 15.243852644089003
$ws.Range("G9").Value = 0.53747886252174237
$ws.Range("H9").Value = 270
$ws.Range("I9").Value = 236
$ws.Range("A10").Value = 18.333333333333336
$ws.Range("B10").Value = 32.192039123608993
$ws.Range("C10").Value = 16.288651714949918
$ws.Range("D10").Value = 0.50598384440344912
$ws.Range("E10").Value = 29.640841367943622
$ws.Range("F10").Value = 15.489875233119628
$ws.Range("G10").Value = 0.52258554475015095
$ws.Range("H10").Value = 276
$ws.Range("I10").Value = 242
$ws.Range("A11").Value = 20
$ws.Range("B11").Value = 33.39462924544231
$ws.Range("C11").Value = 16.317367183436009
$ws.Range("D11").Value = 0.48862249865112661
$ws.Range("E11").Value = 30.590676489596902
$ws.Range("F11").Value = 15.480958121910632
$ws.Range("G11").Value = 0.50606785787085506
$ws.Range("H11").Value = 281.5
$ws.Range("I11").Value = 247
